$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 956.5
$ws.Range("I9").Value = 614
$ws.Range("J9").Value = 1299
$ws.Range("K9").Value = 614
$ws.Range("L9").Value = 1299
$ws.Range("M9").Value = -445
$ws.Range("N9").Value = -1637
$ws.Range("H17").Value = 1810.2572
$ws.Range("I17").Value = 1285.7142
$ws.Range("J17").Value = 2159.9524
$ws.Range("K17").Value = 3857.1426
$ws.Range("L17").Value = 6479.8572
$ws.Range("M17").Value = -3689.1426
$ws.Range("N17").Value = -6815.8572
$ws.Range("H53").Value = 592.8333
$ws.Range("I53").Value = 635.44446
$ws.Range("J53").Value = 465
$ws.Range("K53").Value = 635.44446
$ws.Range("L53").Value = 465
$ws.Range("M53").Value = 1.555539999999951
$ws.Range("N53").Value = -1739
$ws.Range("H74").Value = 12456.23
$ws.Range("I74").Value = 3655.1667
$ws.Range("K74").Value = 3655.1667
$ws.Range("M74").Value = -2719.1667
$ws.Range("H77").Value = 12456.23
$ws.Range("I77").Value = 3655.1667
$ws.Range("K77").Value = 18275.8335
$ws.Range("M77").Value = -13595.8335
$ws.Range("H92").Value = 1123.8
$ws.Range("J92").Value = 2227.5715
$ws.Range("L92").Value = 2227.5715
$ws.Range("N92").Value = -4723.5715
$ws.Range("H96").Value = 418.33334
$ws.Range("I96").Value = 302.14285
$ws.Range("K96").Value = 906.4285500000001
$ws.Range("M96").Value = 466.5714499999999
$ws.Range("H100").Value = 3495.5833
$ws.Range("I100").Value = 1243.375
$ws.Range("K100").Value = 1243.375
$ws.Range("M100").Value = -702.375
$ws.Range("H132").Value = 19148.285
$ws.Range("I132").Value = 26534.625
$ws.Range("K132").Value = 79603.875
$ws.Range("M132").Value = -77073.875
$ws.Range("H138").Value = 4161.2
$ws.Range("J138").Value = 7136
$ws.Range("L138").Value = 21408
$ws.Range("N138").Value = -31688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2575.9333
$ws.Range("I2").Value = 2777.182
$ws.Range("J2").Value = 2022.5
$ws.Range("K2").Value = 2777.182
$ws.Range("L2").Value = 2022.5
$ws.Range("M2").Value = -2664.182
$ws.Range("N2").Value = -2248.5
$ws.Range("H32").Value = 3454.3125
$ws.Range("I32").Value = 3454.3125
$ws.Range("K32").Value = 3454.3125
$ws.Range("M32").Value = -3167.3125
$ws.Range("H34").Value = 66997.5
$ws.Range("J34").Value = 66997.5
$ws.Range("L34").Value = 66997.5
$ws.Range("N34").Value = -67539.5
$ws.Range("H45").Value = 3359.7273
$ws.Range("I45").Value = 2565.2856
$ws.Range("J45").Value = 4750
$ws.Range("K45").Value = 2565.2856
$ws.Range("L45").Value = 4750
$ws.Range("M45").Value = -2188.2856
$ws.Range("N45").Value = -5504
$ws.Range("H61").Value = 5992.1875
$ws.Range("I61").Value = 5716
$ws.Range("K61").Value = 5716
$ws.Range("M61").Value = -5504
$ws.Range("H63").Value = 2924.7144
$ws.Range("I63").Value = 2416.3333
$ws.Range("J63").Value = 3306
$ws.Range("K63").Value = 2416.3333
$ws.Range("L63").Value = 3306
$ws.Range("M63").Value = -1730.3333
$ws.Range("N63").Value = -4678
$ws.Range("H66").Value = 2924.7144
$ws.Range("I66").Value = 2416.3333
$ws.Range("J66").Value = 3306
$ws.Range("K66").Value = 12081.6665
$ws.Range("L66").Value = 16530
$ws.Range("M66").Value = -8649.666499999999
$ws.Range("N66").Value = -23394
$ws.Range("H116").Value = 2575.9333
$ws.Range("I116").Value = 2777.182
$ws.Range("J116").Value = 2022.5
$ws.Range("K116").Value = 2777.182
$ws.Range("L116").Value = 2022.5
$ws.Range("M116").Value = -483.1819999999998
$ws.Range("N116").Value = -6610.5
$ws.Range("H122").Value = 1838.8
$ws.Range("I122").Value = 1838.8
$ws.Range("K122").Value = 5516.4
$ws.Range("M122").Value = -3066.4
$ws.Range("H132").Value = 3331.5833
$ws.Range("I132").Value = 1732.8572
$ws.Range("K132").Value = 5198.571599999999
$ws.Range("M132").Value = -2668.571599999999
$ws.Range("H136").Value = 5992.1875
$ws.Range("I136").Value = 5716
$ws.Range("K136").Value = 17148
$ws.Range("M136").Value = -14598

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2575.9333
$ws.Range("I3").Value = 2777.182
$ws.Range("J3").Value = 2022.5
$ws.Range("K3").Value = 2777.182
$ws.Range("L3").Value = 2022.5
$ws.Range("M3").Value = -2663.182
$ws.Range("N3").Value = -2250.5
$ws.Range("H20").Value = 7568.857
$ws.Range("I20").Value = 8274.6
$ws.Range("K20").Value = 8274.6
$ws.Range("M20").Value = -8027.6
$ws.Range("H105").Value = 1627.5714
$ws.Range("I105").Value = 1482.1666
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1482.1666
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 264.8334
$ws.Range("N105").Value = -5994
$ws.Range("H107").Value = 5203.737
$ws.Range("I107").Value = 1208
$ws.Range("J107").Value = 8799.9
$ws.Range("K107").Value = 1208
$ws.Range("L107").Value = 8799.9
$ws.Range("M107").Value = 712
$ws.Range("N107").Value = -12639.9
$ws.Range("H134").Value = 3564.1428
$ws.Range("I134").Value = 2896.75
$ws.Range("K134").Value = 8690.25
$ws.Range("M134").Value = -6155.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1360.5
$ws.Range("I105").Value = 1508.8889
$ws.Range("K105").Value = 1508.8889
$ws.Range("M105").Value = 238.1111000000001
$ws.Range("H122").Value = 1468.0588
$ws.Range("I122").Value = 1583.3334
$ws.Range("K122").Value = 4750.0002
$ws.Range("M122").Value = -2300.0002
$ws.Range("H132").Value = 3152.7727
$ws.Range("I132").Value = 2659.5557
$ws.Range("K132").Value = 7978.6671
$ws.Range("M132").Value = -5448.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1234
$ws.Range("I5").Value = 1040.5714
$ws.Range("J5").Value = 1572.5
$ws.Range("K5").Value = 3121.7142
$ws.Range("L5").Value = 4717.5
$ws.Range("M5").Value = -3009.7142
$ws.Range("N5").Value = -4941.5
$ws.Range("H121").Value = 813
$ws.Range("J121").Value = 1355.4
$ws.Range("L121").Value = 4066.2
$ws.Range("N121").Value = -6686.200000000001
$ws.Range("H129").Value = 1212.8334
$ws.Range("I129").Value = 855.4
$ws.Range("K129").Value = 2566.2
$ws.Range("M129").Value = 2433.8
$ws.Range("H131").Value = 1961
$ws.Range("I131").Value = 1853.5
$ws.Range("J131").Value = 2498.5
$ws.Range("K131").Value = 5560.5
$ws.Range("L131").Value = 7495.5
$ws.Range("M131").Value = -520.5
$ws.Range("N131").Value = -17575.5
$ws.Range("H132").Value = 2051.5
$ws.Range("I132").Value = 1591.2858
$ws.Range("J132").Value = 2511.7144
$ws.Range("K132").Value = 14321.5722
$ws.Range("L132").Value = 22605.4296
$ws.Range("M132").Value = -11791.5722
$ws.Range("N132").Value = -27665.4296
$ws.Range("H135").Value = 1234
$ws.Range("I135").Value = 1040.5714
$ws.Range("J135").Value = 1572.5
$ws.Range("K135").Value = 9365.142600000001
$ws.Range("L135").Value = 14152.5
$ws.Range("M135").Value = -6830.142600000001
$ws.Range("N135").Value = -19222.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7151.923
$ws.Range("J68").Value = 8303.714
$ws.Range("L68").Value = 8303.714
$ws.Range("N68").Value = -9801.714
$ws.Range("H71").Value = 7151.923
$ws.Range("J71").Value = 8303.714
$ws.Range("L71").Value = 41518.57
$ws.Range("N71").Value = -49006.57
$ws.Range("H74").Value = 34000
$ws.Range("I74").Value = 34000
$ws.Range("K74").Value = 34000
$ws.Range("M74").Value = -33002
$ws.Range("H77").Value = 34000
$ws.Range("I77").Value = 34000
$ws.Range("K77").Value = 102000
$ws.Range("M77").Value = -97008
$ws.Range("H82").Value = 2876.3333
$ws.Range("J82").Value = 3916.4443
$ws.Range("L82").Value = 3916.4443
$ws.Range("N82").Value = -4638.4443
$ws.Range("H85").Value = 2876.3333
$ws.Range("J85").Value = 3916.4443
$ws.Range("L85").Value = 3916.4443
$ws.Range("N85").Value = -6412.4443
$ws.Range("H93").Value = 1777.7778
$ws.Range("I93").Value = 1777.7778
$ws.Range("K93").Value = 1777.7778
$ws.Range("M93").Value = -529.7778000000001
$ws.Range("H100").Value = 6434.737
$ws.Range("I100").Value = 2907.2856
$ws.Range("J100").Value = 8492.416999999999
$ws.Range("K100").Value = 2907.2856
$ws.Range("L100").Value = 8492.416999999999
$ws.Range("M100").Value = -2366.2856
$ws.Range("N100").Value = -9574.416999999999
$ws.Range("H132").Value = 5717.2856
$ws.Range("I132").Value = 4210.3335
$ws.Range("J132").Value = 6847.5
$ws.Range("K132").Value = 12631.0005
$ws.Range("L132").Value = 20542.5
$ws.Range("M132").Value = -10101.0005
$ws.Range("N132").Value = -25602.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 720
$ws.Range("I96").Value = 1050
$ws.Range("K96").Value = 1050
$ws.Range("M96").Value = 323
